$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sample applicants were only there for testing the layout; clear them
# out so the sheet is ready for real data again. Keep row 2 itself (with
# its date-formatted A2 cell) as the first blank entry row.
$ws.Range("A3:N10").Delete()
$ws.Range("A2:N2").ClearContents()

# "POSITION APPLIED" needed a bit more room.
$ws.Range("C1").ColumnWidth = 28

# Leave the cursor on the blank row, ready for the next applicant.
$ws.Range("A2:L2").Select()
